$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cells A2:A6 with consolidated card data
$ws.Range("A2").Value = "('Memnite', ['{0}', 'Artifact Creature — Construct', '1/1'])"
$ws.Range("A3").Value = "('Memoricide', ['{3}{B}', 'Sorcery', 'Choose a nonland card name. Search target player’s graveyard, hand, and library for any number of cards with that name and exile them. Then that player shuffles their library.'])"
$ws.Range("A4").Value = "('Steel Hellkite', ['{6}', 'Artifact Creature — Dragon', 'Flying', '{2}: Steel Hellkite gets +1/+0 until end of turn.', '{X}: Destroy each nonland permanent with converted mana cost X whose controller was dealt combat damage by Steel Hellkite this turn. Activate this ability only once each turn.', '5/5'])"
$ws.Range("A5").Value = "('Tempered Steel', ['{1}{W}{W}', 'Enchantment', 'Artifact creatures you control get +2/+2.'])"
$ws.Range("A6").Value = "('Wurmcoil Engine', ['{6}', 'Artifact Creature — Wurm', 'Deathtouch, lifelink', 'When Wurmcoil Engine dies, create a 3/3 colorless Wurm artifact creature token with deathtouch and a 3/3 colorless Wurm artifact creature token with lifelink.', '6/6'])"

# Delete rows 7 through 26 (entire rows) which are no longer needed
$ws.Range("A7:A26").EntireRow.Delete()
